# Add season-record columns (Wins / Losses / Ties) to the players table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells AD1:AF1 should look like the rest of the header row (bold,
# centered, thin border) -- copy the formatting from the existing AC1
# header cell instead of re-creating a style, then set the header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-44) gets the same season record.
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = 102
    $ws.Cells.Item($row, 31).Value = 60
    $ws.Cells.Item($row, 32).Value = 0
}
